# Update '想去人数' (interested-count) figures across all sheets
# to match freshly re-scraped totals (gh-pages data refresh at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 8443
$ws1.Range("F7").Value = 2377
$ws1.Range("F13").Value = 1036
$ws1.Range("F14").Value = 1591
$ws1.Range("F15").Value = 2189
$ws1.Range("F16").Value = 45
$ws1.Range("F17").Value = 240
$ws1.Range("F18").Value = 301
$ws1.Range("F19").Value = 2052
$ws1.Range("F21").Value = 1020
$ws1.Range("F22").Value = 803
$ws1.Range("F23").Value = 79
$ws1.Range("F25").Value = 1389
$ws1.Range("F26").Value = 585
$ws1.Range("F27").Value = 1520
$ws1.Range("F29").Value = 275
$ws1.Range("F31").Value = 84
$ws1.Range("F33").Value = 2558

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F27").Value = 29
$ws2.Range("F29").Value = 112

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F11").Value = 205
$ws3.Range("F15").Value = 338
$ws3.Range("F16").Value = 2696
$ws3.Range("F17").Value = 327
$ws3.Range("F18").Value = 164
$ws3.Range("F19").Value = 619

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 8443
$ws4.Range("F8").Value = 205
$ws4.Range("F11").Value = 338
$ws4.Range("F12").Value = 2696
$ws4.Range("F14").Value = 164
$ws4.Range("F15").Value = 1036
$ws4.Range("F16").Value = 1591
$ws4.Range("F18").Value = 619
$ws4.Range("F23").Value = 45
$ws4.Range("F24").Value = 301
$ws4.Range("F27").Value = 1020
$ws4.Range("F28").Value = 79
$ws4.Range("F31").Value = 1389
$ws4.Range("F35").Value = 586
$ws4.Range("F38").Value = 1521
$ws4.Range("F39").Value = 275
$ws4.Range("F40").Value = 29
$ws4.Range("F46").Value = 2558
